# prelude of chart for MOB
# Refresh the forecast-match figures (columns B:L, rows 2-27) on the active
# worksheet with the latest computed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3673
$ws.Range("C2").Value = 51
$ws.Range("D2").Value = 561
$ws.Range("E2").Value = 2482
$ws.Range("F2").Value = 943
$ws.Range("G2").Value = 27
$ws.Range("H2").Value = 14282
$ws.Range("I2").Value = 15925
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = 900
$ws.Range("L2").Value = 16154
$ws.Range("B3").Value = 821
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1.632
$ws.Range("E3").Value = 1.662
$ws.Range("F3").Value = 9
$ws.Range("G3").Value = 20
$ws.Range("H3").Value = 55
$ws.Range("I3").Value = 5738
$ws.Range("J3").Value = 29
$ws.Range("B4").Value = 271
$ws.Range("B5").Value = 803
$ws.Range("D5").Value = 476
$ws.Range("E5").Value = 485
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 4
$ws.Range("H5").Value = 104
$ws.Range("I5").Value = 6224
$ws.Range("J5").Value = 7
$ws.Range("B6").Value = 272
$ws.Range("E6").Value = 24
$ws.Range("F6").Value = 4
$ws.Range("H6").Value = 1667
$ws.Range("B7").Value = 268
$ws.Range("D7").Value = 21
$ws.Range("E7").Value = 21
$ws.Range("I7").Value = 5714
$ws.Range("B8").Value = 775
$ws.Range("D8").Value = 442
$ws.Range("E8").Value = 450
$ws.Range("F8").Value = 4
$ws.Range("G8").Value = 4
$ws.Range("H8").Value = 90
$ws.Range("I8").Value = 8927
$ws.Range("B9").Value = 777
$ws.Range("D9").Value = 68
$ws.Range("E9").Value = 68
$ws.Range("I9").Value = 8971
$ws.Range("B10").Value = 779
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 63
$ws.Range("E10").Value = 70
$ws.Range("F10").Value = 7
$ws.Range("H10").Value = 2273
$ws.Range("I10").Value = 17701
$ws.Range("J10").Value = 2
$ws.Range("B11").Value = 1037
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 390
$ws.Range("E11").Value = 417
$ws.Range("F11").Value = 27
$ws.Range("G11").Value = 3
$ws.Range("H11").Value = 2708
$ws.Range("I11").Value = 26871
$ws.Range("J11").Value = 15
$ws.Range("B12").Value = 858
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 74
$ws.Range("E12").Value = 75
$ws.Range("F12").Value = 1
$ws.Range("H12").Value = 159
$ws.Range("I12").Value = 10000
$ws.Range("J12").Value = 2
$ws.Range("B13").Value = 273
$ws.Range("C13").Value = 8
$ws.Range("D13").Value = 302
$ws.Range("E13").Value = 356
$ws.Range("F13").Value = 45
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = 6372
$ws.Range("I13").Value = 4912
$ws.Range("J13").Value = 7
$ws.Range("B14").Value = 804
$ws.Range("D14").Value = 391
$ws.Range("E14").Value = 396
$ws.Range("F14").Value = 2
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = 51
$ws.Range("I14").Value = 5876
$ws.Range("J14").Value = 7
$ws.Range("B15").Value = 270
$ws.Range("B16").Value = 98
$ws.Range("D16").Value = 387
$ws.Range("E16").Value = 1129
$ws.Range("F16").Value = 462
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = 13859
$ws.Range("I16").Value = 1736
$ws.Range("J16").Value = 7
$ws.Range("K16").Value = 239
$ws.Range("L16").Value = 2479
$ws.Range("B17").Value = 23
$ws.Range("D17").Value = 26
$ws.Range("E17").Value = 44
$ws.Range("F17").Value = 17
$ws.Range("H17").Value = 3864
$ws.Range("B18").Value = 5
$ws.Range("C18").Value = 11
$ws.Range("D18").Value = 376
$ws.Range("E18").Value = 793
$ws.Range("F18").Value = 339
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 4330
$ws.Range("I18").Value = 1176
$ws.Range("J18").Value = 6
$ws.Range("K18").Value = 58
$ws.Range("L18").Value = 731
$ws.Range("B19").Value = 788
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 254
$ws.Range("E19").Value = 278
$ws.Range("F19").Value = 24
$ws.Range("H19").Value = 968
$ws.Range("I19").Value = 4486
$ws.Range("J19").Value = 1
$ws.Range("B20").Value = 1405
$ws.Range("B21").Value = 502
$ws.Range("D21").Value = 70
$ws.Range("E21").Value = 95
$ws.Range("F21").Value = 25
$ws.Range("H21").Value = 8404
$ws.Range("I21").Value = 25094
$ws.Range("J21").Value = 1
$ws.Range("B22").Value = 194
$ws.Range("E22").Value = 4
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 50
$ws.Range("B23").Value = 813
$ws.Range("D23").Value = 23
$ws.Range("E23").Value = 23
$ws.Range("I23").Value = 6957
$ws.Range("B24").Value = 1010
$ws.Range("D24").Value = 58
$ws.Range("E24").Value = 64
$ws.Range("F24").Value = 6
$ws.Range("H24").Value = 2857
$ws.Range("I24").Value = 10047
$ws.Range("B25").Value = 24
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 65
$ws.Range("E25").Value = 95
$ws.Range("F25").Value = 10
$ws.Range("H25").Value = 1053
$ws.Range("I25").Value = 1846
$ws.Range("J25").Value = 2
$ws.Range("K25").Value = 19
$ws.Range("L25").Value = 20
$ws.Range("B26").Value = 140
$ws.Range("D26").Value = 78
$ws.Range("E26").Value = 84
$ws.Range("F26").Value = 6
$ws.Range("H26").Value = 714
$ws.Range("I26").Value = 6184
$ws.Range("J26").Value = 2
$ws.Range("B27").Value = 6
$ws.Range("C27").Value = 11
$ws.Range("D27").Value = 107
$ws.Range("E27").Value = 211
$ws.Range("F27").Value = 92
$ws.Range("H27").Value = 4381
$ws.Range("I27").Value = 561
